$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tiny floating point re-computation tweaks on existing rows (6,7,8,12,13) ---
$ws.Cells.Item(6,5).Value  = 9372.060606060608
$ws.Cells.Item(6,6).Value  = 727.3321970634984

$ws.Cells.Item(7,6).Value  = 618.2012296417555

$ws.Cells.Item(8,6).Value  = 1331.618216449506

$ws.Cells.Item(12,6).Value = 976.2876802085996

$ws.Cells.Item(13,3).Value = 12995.45454545455
$ws.Cells.Item(13,4).Value = 35900
$ws.Cells.Item(13,5).Value = 23051.75757575758
$ws.Cells.Item(13,6).Value = 3678.743988305898

# --- Append new rows 14-20 (Season 3 data) ---
# Pull formatting from row 13 (the last existing data row) so the new rows match
# the existing look (bold/centered/bordered season number in column A) without
# introducing brand new style entries.

$newRows = @(
    @{ Row=14; A=12; B="M3_01 Wolf 2021";    C=7968.181818181818;  D=13700; E=9929.636363636364;  F=1046.105394244657 },
    @{ Row=15; A=13; B="M3_02 Love 2021";    C=11645.45454545455;  D=19500; E=13977.75757575757;  F=1430.087790687181 },
    @{ Row=16; A=14; B="M3_03 Bear 2021";    C=7466.666666666667;  D=10900; E=9195.333333333332;  F=715.137593598447  },
    @{ Row=17; A=15; B="M3_04 Elf 2021";     C=11586.36363636364;  D=20300; E=14814.60606060606;  F=1527.809484803296 },
    @{ Row=18; A=16; B="M3_05 Viper 2021";   C=12990.90909090909;  D=20900; E=16159.51515151515;  F=1322.54767637699  },
    @{ Row=19; A=17; B="M3_06 Magic 2021";   C=11263.63636363636;  D=17200; E=13526.06060606061;  F=1034.558786372709 },
    @{ Row=20; A=18; B="M3_07 Griffin 2021"; C=10150;              D=15800; E=12255.33333333333;  F=1040.526362525771 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $prevRow = $rowNum - 1

    # Copy formats only from the row above, so the new row inherits the same
    # per-column styling (e.g. column A's bold/bordered/centered style) without
    # creating extra unused style entries in the workbook.
    $ws.Range("A$prevRow`:F$prevRow").Copy()
    $ws.Range("A$rowNum`:F$rowNum").PasteSpecial(-4122)

    $ws.Cells.Item($rowNum,1).Value = $r.A
    $ws.Cells.Item($rowNum,2).Value = $r.B
    $ws.Cells.Item($rowNum,3).Value = $r.C
    $ws.Cells.Item($rowNum,4).Value = $r.D
    $ws.Cells.Item($rowNum,5).Value = $r.E
    $ws.Cells.Item($rowNum,6).Value = $r.F
}
